# Restore the "From" value of rule R30 (row 10) on the Rules sheet
# from 18 to 1, matching revision #0ce084e2a09761ca63fc5e3164f09828a3eb1340.TEST
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1
